{"js": "// Replace the 25 division-problem answers in the single 20x5 table.\n// The table has data only in rows 0, 4, 8, 12, 16 (every 4th row is a\n// \"work area\" followed by 3 blank rows), 5 answers per data row, for a\n// total of 25 text cells that need to change, in row-major document order.\n\nconst newValues = [\n  [\"25\u00f72=12, 1\", \"68\u00f79=7, 5\", \"27\u00f72=13, 1\", \"90\u00f72=45, 0\", \"89\u00f76=14, 5\"],\n  [\"99\u00f73=33, 0\", \"44\u00f79=4, 8\", \"33\u00f76=5, 3\", \"30\u00f75=6, 0\", \"25\u00f77=3, 4\"],\n  [\"42\u00f75=8, 2\", \"90\u00f72=45, 0\", \"88\u00f75=17, 3\", \"23\u00f76=3, 5\", \"92\u00f73=30, 2\"],\n  [\"42\u00f78=5, 2\", \"70\u00f79=7, 7\", \"90\u00f74=22, 2\", \"20\u00f72=10, 0\", \"87\u00f79=9, 6\"],\n  [\"89\u00f76=14, 5\", \"82\u00f74=20, 2\", \"57\u00f74=14, 1\", \"86\u00f72=43, 0\", \"50\u00f75=10, 0\"]\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nlet dataRowIdx = 0;\nfor (let r = 0; r < table.rowCount; r++) {\n  const rowValues = table.values[r];\n  // A \"data\" row is one that actually has text in its first cell;\n  // the 3 blank spacer rows following each data row are left untouched.\n  const hasData = rowValues && rowValues.some((v) => v && v.trim().length > 0);\n  if (!hasData) continue;\n\n  const replacements = newValues[dataRowIdx];\n  for (let c = 0; c < rowValues.length; c++) {\n    table.getCell(r, c).value = replacements[c];\n  }\n  dataRowIdx++;\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 division-problem answers in the single 20x5 table.\n# The table has data only in rows 1, 5, 9, 13, 17 (1-based; every 4th row\n# is a \"work area\" followed by 3 blank rows), 5 answers per data row, for\n# a total of 25 text cells that need to change, in row-major document\n# order.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newValues = @(\n    @(\"25\u00f72=12, 1\", \"68\u00f79=7, 5\", \"27\u00f72=13, 1\", \"90\u00f72=45, 0\", \"89\u00f76=14, 5\"),\n    @(\"99\u00f73=33, 0\", \"44\u00f79=4, 8\", \"33\u00f76=5, 3\", \"30\u00f75=6, 0\", \"25\u00f77=3, 4\"),\n    @(\"42\u00f75=8, 2\", \"90\u00f72=45, 0\", \"88\u00f75=17, 3\", \"23\u00f76=3, 5\", \"92\u00f73=30, 2\"),\n    @(\"42\u00f78=5, 2\", \"70\u00f79=7, 7\", \"90\u00f74=22, 2\", \"20\u00f72=10, 0\", \"87\u00f79=9, 6\"),\n    @(\"89\u00f76=14, 5\", \"82\u00f74=20, 2\", \"57\u00f74=14, 1\", \"86\u00f72=43, 0\", \"50\u00f75=10, 0\")\n)\n\n$dataRowIdx = 0\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    $firstCellText = $t.Cell($r, 1).Range.Text -replace \"[\\r\\a]\", \"\"\n    if ($firstCellText.Length -eq 0) {\n        continue\n    }\n\n    $replacements = $newValues[$dataRowIdx]\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $t.Cell($r, $c).Range.Text = $replacements[$c - 1]\n    }\n    $dataRowIdx = $dataRowIdx + 1\n}\n"}
